$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Air Separation")
$ws2 = $wb.Worksheets.Item("Lime Kiln")

# ---------------------------------------------------------------
# "Lime Kiln" sheet: rechecking of emission / energy data.
# Insert two new rows (4 & 5) splitting the old single CO2 outflow
# of the calcination step into an intermediate "CO2__temp" quantity
# that is then divided (via Ratio / Remainder) into CO2 retained in
# the slag ("CO2__embodied in lime/slag") and CO2 actually emitted
# ("CO2__emitted"), using a new "slag" known-quantity.
# ---------------------------------------------------------------

$ws2.Rows.Item(4).Insert()
$ws2.Rows.Item(4).Insert()

# New row 4: slag | CO2__temp | temp | CO2__embodied in lime/slag | outflows | Ratio | CO2 not calcinated
$ws2.Cells.Item(4, 1).Value = "slag"
$ws2.Cells.Item(4, 2).Value = "CO2__temp"
$ws2.Cells.Item(4, 3).Value = "temp"
$ws2.Cells.Item(4, 4).Value = "CO2__embodied in lime/slag"
$ws2.Cells.Item(4, 5).Value = "outflows"
$ws2.Cells.Item(4, 6).Value = "Ratio"

# New row 5: slag | CO2__temp | temp | CO2__emitted | outflows | Remainder | CO2 not calcinated
$ws2.Cells.Item(5, 1).Value = "slag"
$ws2.Cells.Item(5, 2).Value = "CO2__temp"
$ws2.Cells.Item(5, 3).Value = "temp"
$ws2.Cells.Item(5, 4).Value = "CO2__emitted"
$ws2.Cells.Item(5, 5).Value = "outflows"
$ws2.Cells.Item(5, 6).Value = "Remainder"

# Variable labels in column G carry the default (unstyled) font.
$ws2.Cells.Item(4, 7).Value = "CO2 not calcinated"
$ws2.Cells.Item(4, 7).Style = "Normal"
$ws2.Cells.Item(5, 7).Value = "CO2 not calcinated"
$ws2.Cells.Item(5, 7).Style = "Normal"

$ws2.Rows.Item(4).RowHeight = 16
$ws2.Rows.Item(5).RowHeight = 16

# Existing row 3 (CaCO3 -> CO2) now feeds the temporary CO2 quantity
# instead of being an outflow directly.
$ws2.Cells.Item(3, 4).Value = "CO2__temp"
$ws2.Cells.Item(3, 5).Value = "temp"

# ---------------------------------------------------------------
# View/selection state: "Lime Kiln" becomes the active tab, with a
# new selection; "Air Separation" keeps its own (different) selection.
# ---------------------------------------------------------------
$ws1.Range("B8").Select()
$ws2.Activate()
$ws2.Range("D7").Select()
